$wb = $excel.ActiveWorkbook

# Insert two new columns (B, C) on every sheet, shifting the old p.value column to D.
# Set the new headers and fill in coefficient.estimates / standard.deviation values.

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("B:C").Insert()
$ws.Range("B1").Value = "coefficient.estimates"
$ws.Range("C1").Value = "standard.deviation"

$coef1 = @(
-0.0763052271519631,0.10033794990673,-0.0670021225634853,-0.0705368429346508,-0.0604666234345592,-0.0561617628829433,-0.0554653492159775,-0.148935707213245,-0.084668096032864,-0.0779037691436003,-0.0584013075596875,-0.0420657040715558,-0.0494180754129195,-0.0486824963578774,-0.065561781402192,-0.0395532755819045,-0.0721404039918988,-0.0509741272935522,-0.049483371211158,-0.0444286046878637,-0.048968729734294,-0.0543378566316912,-0.0432736133160775,-0.0364013002509875,-0.0349417210926487,-0.0328797587547415,-0.0212038487805496,-0.041883650622317,-0.0500639422743185,-0.0355943380605704,-0.0252079200497155,-0.0279185664179544,-0.0278728594887648,-0.0367090623983753,-0.0353908986131646,-0.0339265137663662,-0.0516418237771288,-0.0461392575220597,-0.0173510820931784,-0.0385285867550733,-0.0238794242618883,-0.0258276465498037,-0.0274843857977052,-0.0337216484307791,-0.031618543711765,-0.0194091004097648,-0.0428132488903195,-0.0364650457578673,0.0156003580378038,-0.0284209650370402,-0.0222945556673554,0.0146636746639555,-0.0400598306576633,-0.0323433640209836,-0.0190941353475554,-0.0117897634023689,-0.0237894073887506,0.00676424536639544,-0.0040742252784196,0.00372214670965654,-0.00565383103112839
)
$sd1 = @(
0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515,0.0374806734950515
)
for ($i = 0; $i -lt 61; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $coef1[$i]
    $ws.Cells.Item($row, 3).Value = $sd1[$i]
}

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("B:C").Insert()
$ws.Range("B1").Value = "coefficient.estimates"
$ws.Range("C1").Value = "standard.deviation"

$coef2 = @(
-0.258951327648078,-0.394398268721543,-0.353762341210877,-0.407976937534661,-0.299476972190242,0.166370887077954,-0.323522320641329,-0.252068676596143,0.40512488191822,-0.266773763805851,-0.322439953646661,0.227915225389179,-0.0891331238469593,-0.348857674419955,0.332476097811617,0.23250983397156,-0.219793615792888,-0.145380446791152,-0.176461165068819,-0.205650339394867,-0.0716287252628522,-0.104730299462865,-0.158873646200483,-0.191414757904197,0.0882901234465362,0.110905842805974,0.178113542974048,-0.190780357189742,-0.0922909338704116,0.147429657314894,-0.173667074513425,-0.205750321315987,0.091591350521328,-0.0950040387726264,-0.216481745328884,-0.115402564186821,0.103448677494892,-0.133017398665326,0.0193916498001613,-0.0354198793170126,-0.131675055306126,0.126612137391757,0.0161208588778511,-0.092252104403152,0.0536124888081878,0.0529423486123759,0.0788369889026203,0.0543680766162146,0.0642888221274042,-0.0862360054173692,-0.13152900527624,-0.0871025162964395,-0.0444307347738558,-0.0234878068779793,-0.0136012128146742,0.0323625090679225,0.0358609800725144,0.00991973204926727,-0.0300921661170291,0.000105190053618457,0.0000723447881021558
)
$sd2 = @(
0.100955679803525,0.171708786931031,0.156304104574164,0.182350931194948,0.134787617936956,0.0750541478737976,0.149604850127314,0.118505654474348,0.197963131939006,0.1322255102074,0.163194508685278,0.122076556557597,0.0479278109921588,0.188339920615496,0.203630010100573,0.156927743026904,0.152204134938873,0.100919070477466,0.147198883847461,0.187235183205428,0.0660202589945298,0.0965601382117368,0.152317720493298,0.183682665284889,0.0868149389200065,0.111435941019205,0.183192650939163,0.196931197164489,0.0960964304780509,0.155556623413392,0.189030827437811,0.236026302146126,0.105498840882915,0.111318114409138,0.257608852588407,0.14446456067878,0.134041775603845,0.17337487284305,0.0263321423814016,0.051483885055862,0.195368488894445,0.190838926652202,0.0250311902733534,0.145499852133536,0.0963632987130167,0.0964705426751755,0.150352242104167,0.106941895911069,0.135300146401151,0.181641552886474,0.292789237302091,0.196005949994992,0.141087036842489,0.0754250255723894,0.0446534036570745,0.114068422580879,0.141425798016166,0.0434644490094733,0.163240953876681,0.193020156850397,0.140024572618724
)
for ($i = 0; $i -lt 61; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $coef2[$i]
    $ws.Cells.Item($row, 3).Value = $sd2[$i]
}

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("B:C").Insert()
$ws.Range("B1").Value = "coefficient.estimates"
$ws.Range("C1").Value = "standard.deviation"

$coef3 = @(
-0.246017579632504,-0.564231201667572,-0.405661799476967,-0.33299494489613,0.177533173480202,-0.334264837877482,-0.22821567867805,-0.294522505565298,0.383792434194368,0.287126544916958,-0.0949043390664925,-0.281895140376555,-0.218913246130673,0.460613497157691,0.367983926218481,-0.312158168785352,-0.171721043969404,0.242102937248176,0.231441135118057,0.249797351429648,-0.227180715764979,-0.192401791975094,0.243454142454398,0.190801639780719,-0.266049508667092,-0.250209035684544,-0.23201846844974,-0.232069976930907,0.226542202962918,-0.0731552492123801,-0.116412720276773,-0.0867477118019402,0.0738329507461105,-0.113011988700372,-0.172357781441202,0.165192512826169,-0.182967120835435,-0.0446291211149562,-0.126981379063558,0.113253202488209,-0.189897007350177,0.0607107396470452,0.0177223896354788,0.0168932100622719,0.127373885631913,-0.132833449775561,-0.0849811266740541,-0.105204489166354,-0.147311074139622,-0.0238979035541535,-0.0421875922536895,0.0316541562231578,-0.0507297473477109,0.0237069516679031,-0.0519430499631489,-0.00637492085386211,-0.0329602821395344,-0.0307801763809465,0.0107605859463177,0.00415355524487923,0.00915363957186925
)
$sd3 = @(
0.0983720523565027,0.227809333500032,0.171349505317142,0.142543121566335,0.0779409922663357,0.155728482620262,0.111144183710231,0.146519711969838,0.194327825650077,0.147680278420714,0.0508698758903222,0.151179945820018,0.119552625054419,0.2515735659266,0.206302119772561,0.178740773685463,0.0992654115350651,0.141530366958767,0.13713700845308,0.167426944105893,0.154342533070782,0.132568037897244,0.194137457740883,0.156640348969051,0.218997185661163,0.207583521719713,0.199755672250185,0.205077950625962,0.207537448470885,0.0674331114430182,0.109833447859149,0.0858453231994032,0.0816161407774001,0.125303702348484,0.191605218112843,0.195975245998127,0.217430918729944,0.0571295081767383,0.167667708953629,0.151518775342546,0.256081226711238,0.0857246028882516,0.0271171427418133,0.0259333269676264,0.205867034250684,0.215417476144703,0.147426649949108,0.224688836913596,0.35827249868693,0.0597178726313914,0.136963902084974,0.104632378547285,0.19479807469331,0.10192556922259,0.258250441193387,0.035526871043902,0.199966766570015,0.201078990081861,0.102797455509316,0.040883999452975,0.151453746655061
)
for ($i = 0; $i -lt 61; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $coef3[$i]
    $ws.Cells.Item($row, 3).Value = $sd3[$i]
}

# ---- Sheet 4 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("B:C").Insert()
$ws.Range("B1").Value = "coefficient.estimates"
$ws.Range("C1").Value = "standard.deviation"

$coef4 = @(
0.515028408149602,0.44659940987466,0.19045377129479,0.392582550531874,0.366932300247328,0.30183384056,0.258714704159846,0.364869356523869,0.370579579047325,0.322594513910633,0.635776667747003,0.334681691582255,0.344503902104763,0.342413482293668,0.156082724635449,0.162884812503849,0.256880537407976,0.357360676582343,-0.149917603745626,0.0865130919710721,0.128870791246848,-0.23049188841008,0.21277651967362,-0.056215191210046,0.105877890605968,0.163800662871134,0.203550186945292,0.248792266526209,0.123299825175395,0.204741836140257,0.252019310550381,0.101804247212524,0.223270125919901,0.149468285519394,-0.171699475055042,0.21143440222555,0.189563930943986,-0.0766639768694214,-0.104981009798922,-0.0851499563702806,0.0128802079294336,0.0124309642088046,0.0992838384510815,-0.11769330524445,-0.0778000840590078,-0.101754953534454,-0.024458509649832,0.123806872442572,-0.0505041396806549,-0.0559562264610314,0.037521606794193,-0.0361662735780244,0.00850552877815013,0.0229231970001317,-0.0264931359947918,-0.022483130941447,0.0145383841171041,0.0053057231841713,-0.00391512191909448,-0.000159420773231722
)
$sd4 = @(
0.157305365749482,0.138159070737969,0.0589351961955329,0.123407746501872,0.116149200281331,0.104752773523864,0.0936690095821512,0.134738035845257,0.14014689997368,0.128670571119148,0.253874803288348,0.139291605589631,0.163565448923016,0.181659910760811,0.0889420303703149,0.0945388643773125,0.156396227241,0.218655991916893,0.0937097177518782,0.0564928389275608,0.0936718561566443,0.174660450055332,0.162556783429128,0.0434376454945668,0.0822662429779459,0.12732259807492,0.161415170342574,0.202080955672493,0.111459478121117,0.186926571503232,0.238788015234305,0.0989842331602398,0.220347471798054,0.147855379892862,0.179012849810941,0.221205766127167,0.20217973905644,0.088780747142355,0.126078553826087,0.106812312946574,0.0161903027998924,0.015644142756302,0.126113986089579,0.149560476835847,0.102862258913295,0.137912510108687,0.0401949001773027,0.216015782688411,0.112715798909774,0.126422657661321,0.101651678015525,0.11298171077179,0.0275959380071587,0.146120510774828,0.187782161089443,0.167025302120119,0.120240503461609,0.0522366281179571,0.0493211537681793,0.0398556312967697
)
for ($i = 0; $i -lt 60; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $coef4[$i]
    $ws.Cells.Item($row, 3).Value = $sd4[$i]
}
